$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-01-19 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-01-20 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("501÷9=55, 6", $true, $false, $false, $false, $false, $true, 1, $false, "903÷5=180, 3", 2) | Out-Null
$d.Content.Find.Execute("597÷2=298, 1", $true, $false, $false, $false, $false, $true, 1, $false, "977÷5=195, 2", 2) | Out-Null
$d.Content.Find.Execute("745÷8=93, 1", $true, $false, $false, $false, $false, $true, 1, $false, "684÷8=85, 4", 2) | Out-Null
$d.Content.Find.Execute("340÷8=42, 4", $true, $false, $false, $false, $false, $true, 1, $false, "497÷2=248, 1", 2) | Out-Null
$d.Content.Find.Execute("318÷5=63, 3", $true, $false, $false, $false, $false, $true, 1, $false, "659÷6=109, 5", 2) | Out-Null
$d.Content.Find.Execute("814÷6=135, 4", $true, $false, $false, $false, $false, $true, 1, $false, "265÷4=66, 1", 2) | Out-Null
$d.Content.Find.Execute("848÷7=121, 1", $true, $false, $false, $false, $false, $true, 1, $false, "471÷6=78, 3", 2) | Out-Null
$d.Content.Find.Execute("113÷8=14, 1", $true, $false, $false, $false, $false, $true, 1, $false, "150÷3=50, 0", 2) | Out-Null
$d.Content.Find.Execute("995÷3=331, 2", $true, $false, $false, $false, $false, $true, 1, $false, "821÷2=410, 1", 2) | Out-Null
$d.Content.Find.Execute("583÷4=145, 3", $true, $false, $false, $false, $false, $true, 1, $false, "787÷4=196, 3", 2) | Out-Null
$d.Content.Find.Execute("880÷3=293, 1", $true, $false, $false, $false, $false, $true, 1, $false, "412÷2=206, 0", 2) | Out-Null
$d.Content.Find.Execute("288÷4=72, 0", $true, $false, $false, $false, $false, $true, 1, $false, "466÷3=155, 1", 2) | Out-Null
$d.Content.Find.Execute("275÷5=55, 0", $true, $false, $false, $false, $false, $true, 1, $false, "589÷9=65, 4", 2) | Out-Null
$d.Content.Find.Execute("397÷5=79, 2", $true, $false, $false, $false, $false, $true, 1, $false, "785÷9=87, 2", 2) | Out-Null
$d.Content.Find.Execute("620÷4=155, 0", $true, $false, $false, $false, $false, $true, 1, $false, "631÷3=210, 1", 2) | Out-Null
$d.Content.Find.Execute("802÷7=114, 4", $true, $false, $false, $false, $false, $true, 1, $false, "742÷6=123, 4", 2) | Out-Null
$d.Content.Find.Execute("582÷5=116, 2", $true, $false, $false, $false, $false, $true, 1, $false, "994÷5=198, 4", 2) | Out-Null
$d.Content.Find.Execute("431÷9=47, 8", $true, $false, $false, $false, $false, $true, 1, $false, "357÷2=178, 1", 2) | Out-Null
$d.Content.Find.Execute("782÷6=130, 2", $true, $false, $false, $false, $false, $true, 1, $false, "191÷6=31, 5", 2) | Out-Null
$d.Content.Find.Execute("761÷8=95, 1", $true, $false, $false, $false, $false, $true, 1, $false, "759÷8=94, 7", 2) | Out-Null
$d.Content.Find.Execute("667÷5=133, 2", $true, $false, $false, $false, $false, $true, 1, $false, "282÷7=40, 2", 2) | Out-Null
$d.Content.Find.Execute("229÷6=38, 1", $true, $false, $false, $false, $false, $true, 1, $false, "280÷2=140, 0", 2) | Out-Null
$d.Content.Find.Execute("910÷5=182, 0", $true, $false, $false, $false, $false, $true, 1, $false, "453÷6=75, 3", 2) | Out-Null
$d.Content.Find.Execute("666÷4=166, 2", $true, $false, $false, $false, $false, $true, 1, $false, "218÷9=24, 2", 2) | Out-Null
$d.Content.Find.Execute("152÷4=38, 0", $true, $false, $false, $false, $false, $true, 1, $false, "492÷9=54, 6", 2) | Out-Null
